{"js": "// tdf#108714 follow-up fixture edit:\n// - \"Paragraph 2\" -> \"Paragraph 3\"\n// - \"Paragraph 3\" -> \"Paragraph 4\"\n// - add a one-cell table containing \"Paragraph 5 in table\" right after it\n// - add a page break after the table\n// - add a new paragraph \"Paragraph 6\" after that break\n\nconst body = context.document.body;\n\n// Rename the two existing trailing paragraphs in place.\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"text\");\nawait context.sync();\n\nlet para2 = null;\nlet para3 = null;\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  const t = paragraphs.items[i].text;\n  if (t === \"Paragraph 2\") para2 = paragraphs.items[i];\n  else if (t === \"Paragraph 3\") para3 = paragraphs.items[i];\n}\nif (para2) para2.insertText(\"Paragraph 3\", \"Replace\");\nif (para3) para3.insertText(\"Paragraph 4\", \"Replace\");\nawait context.sync();\n\n// Append the new table holding \"Paragraph 5 in table\" using a minimal OOXML\n// fragment so no extra table/cell formatting gets invented.\nconst tableOoxml = `<?xml version=\"1.0\" encoding=\"UTF-8\" standalone=\"yes\"?>\n<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">\n  <pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\">\n    <pkg:xmlData>\n      <w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\">\n        <w:body>\n          <w:tbl>\n            <w:tr>\n              <w:tc>\n                <w:p>\n                  <w:r>\n                    <w:t xml:space=\"preserve\">Paragraph 5 in table</w:t>\n                  </w:r>\n                </w:p>\n              </w:tc>\n            </w:tr>\n          </w:tbl>\n        </w:body>\n      </w:document>\n    </pkg:xmlData>\n  </pkg:part>\n</pkg:package>`;\nbody.insertOoxml(tableOoxml, \"End\");\nawait context.sync();\n\n// Page break right after the table, then the new trailing paragraph.\nbody.insertBreak(\"Page\", \"End\");\nawait context.sync();\n\nbody.insertParagraph(\"Paragraph 6\", \"End\");\nawait context.sync();\n", "ps1": "# tdf#108714 follow-up fixture edit:\n# - \"Paragraph 2\" -> \"Paragraph 3\"\n# - \"Paragraph 3\" -> \"Paragraph 4\"\n# - add a one-cell table containing \"Paragraph 5 in table\" right after it\n# - add a page break after the table\n# - add a new paragraph \"Paragraph 6\" after that break\n\n$d = $word.ActiveDocument\n\n# Step 1: rename the two trailing paragraphs in place (match on text,\n# ignoring the paragraph-mark CR that Range.Text includes).\nfor ($i = 1; $i -le $d.Paragraphs.Count; $i++) {\n    $p = $d.Paragraphs.Item($i)\n    $t = $p.Range.Text.TrimEnd([char]13)\n    if ($t -eq \"Paragraph 2\") {\n        $p.Range.Text = \"Paragraph 3\"\n    } elseif ($t -eq \"Paragraph 3\") {\n        $p.Range.Text = \"Paragraph 4\"\n    }\n}\n\n# Step 2: append the new table holding \"Paragraph 5 in table\" using a\n# minimal OOXML fragment, so no extra table/cell formatting gets invented.\n$endRange = $d.Range()\n$endRange.Collapse(0)\n$tableXml = '<?xml version=\"1.0\" encoding=\"UTF-8\" standalone=\"yes\"?>\n<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">\n  <pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\">\n    <pkg:xmlData>\n      <w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\">\n        <w:body>\n          <w:tbl>\n            <w:tr>\n              <w:tc>\n                <w:p>\n                  <w:r>\n                    <w:t xml:space=\"preserve\">Paragraph 5 in table</w:t>\n                  </w:r>\n                </w:p>\n              </w:tc>\n            </w:tr>\n          </w:tbl>\n        </w:body>\n      </w:document>\n    </pkg:xmlData>\n  </pkg:part>\n</pkg:package>'\n$endRange.InsertXML($tableXml)\n\n# Step 3: page break right after the table (in its own paragraph), then\n# the new trailing paragraph with the text.\n$breakRange = $d.Range()\n$breakRange.Collapse(0)\n$breakRange.InsertBreak(7)  # wdPageBreak\n$breakRange.InsertParagraphAfter()\n\n$tailRange = $d.Range()\n$tailRange.Collapse(0)\n$tailRange.InsertAfter(\"Paragraph 6\")\n"}
